$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Jared Olmos"
$ws.Range("B2").Value = 4.6
$ws.Range("C2").Value = 5.89
$ws.Range("D2").Value = 40.8
$ws.Range("E2").Value = 3.1
$ws.Range("F2").Value = 40.53
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 3.44
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 23.5
$ws.Range("K2").Value = 74.3
$ws.Range("L2").Value = 12.45
$ws.Range("M2").Value = 60.78
$ws.Range("N2").Value = 5.7
$ws.Range("O2").Value = 56
